$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 currently holds the trailing "empty but styled" cell (A6, same
# look as the other column-A header cells A1:A5). Turn it into a real
# "style" / "default" key-value data row, then recreate a fresh trailing
# empty styled cell one row down (A7), matching the sheet's convention.

# Copy the existing header style (bold/orange font) from A1 so the new
# A6/A7 cells reuse the same style index instead of creating a new one.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "default"
